$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that dropped out of this "missing data" sample:
# "RM 232" (old row 26) and "SC 92" (old row 28). Delete the lower-numbered
# row first so the higher row index is still valid when we delete it.
$ws.Rows("28").Delete()
$ws.Rows("26").Delete()

# After the deletions, the remaining rows have shifted up by one or two.
# Apply the per-cell value changes (including newly-missing / newly-restored
# values) so the sheet matches the target "missing_data" sample exactly.
$ws.Range("E19").Value = -6.5

$ws.Range("E21").ClearContents()

$ws.Range("E23").Value = -7

$ws.Range("D26").ClearContents()

$ws.Range("D27").Value = -14.6
$ws.Range("E27").ClearContents()

$ws.Range("D29").ClearContents()

$ws.Range("E33").Value = -10.7
